# edit.ps1
# Commit: "Add data for 2021-12-20"
#
# This advances the workbook's "through" date from December 11 to December 12,
# 2021, and refreshes the monthly carjacking counts for every neighborhood whose
# numbers changed as a result of the newly-added day of data.
#
# One side effect of the new totals: Englewood's running "through" count (5)
# now exceeds North Lawndale's (3), so the two neighborhoods trade places in
# the (count-sorted) table -- row 3 becomes Englewood and row 4 becomes
# North Lawndale, each carrying its own refreshed monthly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook / sheet level text -------------------------------------------------
$wb.Worksheets.Item(1).Name = "Through 2021-12-12"
$ws.Range("B1").Value = "December 2021 (through December 12)"

# --- Row labels: Englewood now outranks North Lawndale, so the two rows swap ----
$ws.Range("A3").Value = "Englewood"
$ws.Range("A4").Value = "North Lawndale"

# --- Per-cell numeric updates ----------------------------------------------------
# Row 3: Englewood (was North Lawndale's slot)
$ws.Range("B3").Value = 5   # was 3
$ws.Range("C3").Value = 8   # was 9
$ws.Range("D3").Value = 7   # was 19
$ws.Range("E3").Value = 2   # was 10
$ws.Range("F3").Value = 3   # was 12
$ws.Range("G3").Value = 6   # was 9
$ws.Range("H3").Value = 7   # was 8
$ws.Range("I3").Value = 7   # was 2
$ws.Range("J3").Value = 3   # was 9
$ws.Range("K3").Value = 2   # was 10
$ws.Range("L3").Value = $null   # was 7
$ws.Range("M3").Value = 4   # was 10
$ws.Range("N3").Value = 2   # was 7
$ws.Range("O3").Value = 7   # was 19
$ws.Range("P3").Value = 8   # was 12
$ws.Range("Q3").Value = 2   # was 9
$ws.Range("R3").Value = 2   # was 5
$ws.Range("S3").Value = 9   # was 25
$ws.Range("T3").Value = 8   # was 5
$ws.Range("U3").Value = 7   # was 3
$ws.Range("V3").Value = 2   # was None
$ws.Range("W3").Value = 4   # was 1
$ws.Range("X3").Value = 3   # was 6
$ws.Range("Y3").Value = 5   # was 4
$ws.Range("Z3").Value = $null   # was 1
$ws.Range("AA3").Value = 3   # was 7
$ws.Range("AB3").Value = 2   # was 1
$ws.Range("AC3").Value = 1   # was 6
$ws.Range("AD3").Value = 3   # was None
$ws.Range("AE3").Value = 2   # was 1
$ws.Range("AF3").Value = 4   # was 3
$ws.Range("AG3").Value = 2   # was 3
$ws.Range("AH3").Value = 1   # was 4
$ws.Range("AI3").Value = $null   # was 3
$ws.Range("AJ3").Value = 2   # was 1
$ws.Range("AK3").Value = 5   # was 1
$ws.Range("AM3").Value = 3   # was 1
$ws.Range("AO3").Value = 2   # was 3
$ws.Range("AP3").Value = 3   # was 5
$ws.Range("AQ3").Value = 2   # was None
$ws.Range("AR3").Value = 1   # was 6
$ws.Range("AT3").Value = 1   # was 2
$ws.Range("AU3").Value = 7   # was 2
$ws.Range("AV3").Value = 1   # was 3
$ws.Range("AX3").Value = $null   # was 2
$ws.Range("AZ3").Value = 5   # was 4
$ws.Range("BB3").Value = 2   # was 6
$ws.Range("BC3").Value = 1   # was 3
$ws.Range("BE3").Value = $null   # was 2
$ws.Range("BF3").Value = 4   # was 6
$ws.Range("BG3").Value = 1   # was 3
$ws.Range("BH3").Value = 3   # was 1
$ws.Range("BI3").Value = $null   # was 4
$ws.Range("BJ3").Value = $null   # was 2
$ws.Range("BK3").Value = 5   # was 2
$ws.Range("BL3").Value = 5   # was 2
$ws.Range("BM3").Value = 4   # was 1
$ws.Range("BN3").Value = 1   # was 2
$ws.Range("BO3").Value = 5   # was 3
$ws.Range("BP3").Value = 4   # was None
$ws.Range("BQ3").Value = 2   # was 4
$ws.Range("BR3").Value = 4   # was None
$ws.Range("BT3").Value = 2   # was 1
$ws.Range("BU3").Value = 5   # was 4
$ws.Range("BV3").Value = 1   # was None
$ws.Range("BW3").Value = 1   # was 2
$ws.Range("BX3").Value = 1   # was 2
$ws.Range("BY3").Value = 5   # was 2
$ws.Range("BZ3").Value = 3   # was None
$ws.Range("CB3").Value = 2   # was 1
$ws.Range("CC3").Value = 1   # was None
$ws.Range("CD3").Value = 2   # was 1

# Row 4: North Lawndale (was Englewood's slot)
$ws.Range("B4").Value = 3   # was 5
$ws.Range("D4").Value = 19   # was 7
$ws.Range("E4").Value = 10   # was 2
$ws.Range("F4").Value = 12   # was 3
$ws.Range("G4").Value = 9   # was 6
$ws.Range("H4").Value = 8   # was 7
$ws.Range("I4").Value = 2   # was 7
$ws.Range("J4").Value = 9   # was 3
$ws.Range("K4").Value = 10   # was 2
$ws.Range("L4").Value = 7   # was None
$ws.Range("M4").Value = 10   # was 4
$ws.Range("N4").Value = 8   # was 2
$ws.Range("O4").Value = 19   # was 7
$ws.Range("P4").Value = 12   # was 8
$ws.Range("Q4").Value = 9   # was 2
$ws.Range("R4").Value = 5   # was 2
$ws.Range("S4").Value = 25   # was 9
$ws.Range("T4").Value = 5   # was 8
$ws.Range("U4").Value = 3   # was 7
$ws.Range("V4").Value = $null   # was 2
$ws.Range("W4").Value = 1   # was 4
$ws.Range("X4").Value = 6   # was 3
$ws.Range("Y4").Value = 4   # was 5
$ws.Range("Z4").Value = 1   # was None
$ws.Range("AA4").Value = 7   # was 3
$ws.Range("AB4").Value = 1   # was 2
$ws.Range("AC4").Value = 6   # was 1
$ws.Range("AD4").Value = $null   # was 3
$ws.Range("AE4").Value = 1   # was 2
$ws.Range("AF4").Value = 3   # was 4
$ws.Range("AG4").Value = 3   # was 2
$ws.Range("AH4").Value = 4   # was 1
$ws.Range("AI4").Value = 3   # was None
$ws.Range("AJ4").Value = 1   # was 2
$ws.Range("AK4").Value = 1   # was 5
$ws.Range("AM4").Value = 1   # was 3
$ws.Range("AO4").Value = 3   # was 2
$ws.Range("AP4").Value = 5   # was 3
$ws.Range("AQ4").Value = $null   # was 2
$ws.Range("AR4").Value = 6   # was 1
$ws.Range("AT4").Value = 2   # was 1
$ws.Range("AU4").Value = 2   # was 7
$ws.Range("AV4").Value = 3   # was 1
$ws.Range("AX4").Value = 2   # was None
$ws.Range("AZ4").Value = 4   # was 5
$ws.Range("BB4").Value = 6   # was 2
$ws.Range("BC4").Value = 3   # was 1
$ws.Range("BE4").Value = 2   # was None
$ws.Range("BF4").Value = 6   # was 4
$ws.Range("BG4").Value = 3   # was 1
$ws.Range("BH4").Value = 1   # was 3
$ws.Range("BI4").Value = 4   # was None
$ws.Range("BJ4").Value = 2   # was None
$ws.Range("BK4").Value = 2   # was 5
$ws.Range("BL4").Value = 2   # was 5
$ws.Range("BM4").Value = 1   # was 4
$ws.Range("BN4").Value = 2   # was 1
$ws.Range("BO4").Value = 3   # was 5
$ws.Range("BP4").Value = $null   # was 4
$ws.Range("BQ4").Value = 4   # was 2
$ws.Range("BR4").Value = $null   # was 4
$ws.Range("BT4").Value = 1   # was 2
$ws.Range("BU4").Value = 4   # was 5
$ws.Range("BW4").Value = 2   # was 1
$ws.Range("BX4").Value = 2   # was 1
$ws.Range("BY4").Value = 2   # was 5
$ws.Range("BZ4").Value = $null   # was 3
$ws.Range("CB4").Value = 1   # was 2
$ws.Range("CC4").Value = $null   # was 1
$ws.Range("CD4").Value = 1   # was 2

# Row 5: West Pullman
$ws.Range("BV5").Value = 3   # was 2

# Row 7: Austin
$ws.Range("AL7").Value = 4   # was 3

# Row 8: Chatham
$ws.Range("BJ8").Value = 1   # was None

# Row 14: Chicago Lawn
$ws.Range("BJ14").Value = 2   # was 1

# Row 15: Washington Heights
$ws.Range("N15").Value = 4   # was 3
$ws.Range("BJ15").Value = 1   # was None

# Row 18: Grand Boulevard
$ws.Range("N18").Value = 1   # was None

# Row 24: South Shore
$ws.Range("N24").Value = 2   # was 1

# Row 30: South Chicago
$ws.Range("B30").Value = 2   # was 1

# Row 38: Auburn Gresham
$ws.Range("AL38").Value = 2   # was 1
$ws.Range("AX38").Value = 5   # was 4

# Row 40: Calumet Heights
$ws.Range("B40").Value = 2   # was 1
$ws.Range("Z40").Value = 1   # was None

# Row 41: Chinatown
$ws.Range("B41").Value = 3   # was 2

# Row 43: Gage Park
$ws.Range("B43").Value = 1   # was None

# Row 54: Rogers Park
$ws.Range("AL54").Value = 1   # was None

# Row 57: Portage Park
$ws.Range("Z57").Value = 2   # was 1

# Row 93: River North
$ws.Range("B93").Value = 3   # was 2
$ws.Range("BJ93").Value = 1   # was None

# Row 97: Streeterville
$ws.Range("N97").Value = 1   # was None
